# Generate Report for Handback
# Refresh the handback status report:
#  - the 12e9f1b5-... file is no longer in sync with en-US
#  - the cf5d4e63-... file has been handed back again (new handoff/handback timestamps)
# Also widen the long "status/date" columns slightly (report column auto-sizing).

$wb = $excel.ActiveWorkbook

$notInSync = "Handed back: not in sync with en-US"

# ----- Overview sheet -----
$wsOverview = $wb.Worksheets.Item("Overview")

# 12e9f1b5 row (row 2) flips from "in sync" to "not in sync"
$wsOverview.Range("E2").Value = $notInSync
$wsOverview.Range("F2").Value = $notInSync

# cf5d4e63 row (row 3) gets a fresh "Latest HO Xliff Generate Date"
$wsOverview.Range("G3").Value = "2017-02-21 11:01:37"

# Widen status columns to fit the new, longer text
$wsOverview.Columns.Item(5).ColumnWidth = 32.6
$wsOverview.Columns.Item(6).ColumnWidth = 32.6

# ----- zh-cn sheet -----
$wsZh = $wb.Worksheets.Item("zh-cn")

# 12e9f1b5 row (row 2) flips from "in sync" to "not in sync"
$wsZh.Range("C2").Value = $notInSync

# cf5d4e63 row (row 3) gets new handoff/handback timestamps
$wsZh.Range("H3").Value = "2017-02-21 11:01:20"
$wsZh.Range("L3").Value = "2017-02-21 11:02:19"

# Widen status column to fit the new, longer text
$wsZh.Columns.Item(3).ColumnWidth = 32.6

# ----- de-de sheet -----
$wsDe = $wb.Worksheets.Item("de-de")

# 12e9f1b5 row (row 2) flips from "in sync" to "not in sync"
$wsDe.Range("C2").Value = $notInSync

# cf5d4e63 row (row 3) gets new handoff/handback timestamps
$wsDe.Range("H3").Value = "2017-02-21 11:01:37"
$wsDe.Range("L3").Value = "2017-02-21 11:02:41"

# Widen status column to fit the new, longer text
$wsDe.Columns.Item(3).ColumnWidth = 32.6
